$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "codeforces_579_A.java"
$ws.Range("D3").Value = "atcoder_ABC169_D.java"
$ws.Range("D4").Value = "atcoder_ARC102_C.java"
$ws.Range("D5").Value = "atcoder_ABC108_B.java"
$ws.Range("D6").Value = "atcoder_ABC120_C.java"
$ws.Range("D7").Value = "codeforces_672_A.java"
$ws.Range("D8").Value = "codeforces_55_A.java"
$ws.Range("D9").Value = "codeforces_110_B.java"
$ws.Range("D10").Value = "codeforces_678_B.java"
$ws.Range("D12").Value = "codeforces_147_A.java"
$ws.Range("D13").Value = "codeforces_189_A.java"
$ws.Range("D14").Value = "codeforces_58_B.java"
$ws.Range("D15").Value = "atcoder_ABC122_D.java"
$ws.Range("D17").Value = "atcoder_ABC132_F.java"
$ws.Range("D18").Value = "codeforces_99_A.java"
$ws.Range("D19").Value = "atcoder_AGC025_A.java"
$ws.Range("D20").Value = "codeforces_92_A.java"
$ws.Range("D21").Value = "codeforces_514_A.java"
$ws.Range("D22").Value = "codeforces_8_B.java"
$ws.Range("D23").Value = "atcoder_ABC178_A.java"
$ws.Range("D24").Value = "atcoder_ABC124_C.java"
$ws.Range("D25").Value = "atcoder_ABC172_D.java"
$ws.Range("D26").Value = "atcoder_ABC136_B.java"
$ws.Range("D27").Value = "codeforces_32_B.java"
$ws.Range("D28").Value = "codeforces_86_A.java"
$ws.Range("D30").Value = "atcoder_ABC114_C.java"
$ws.Range("D31").Value = "atcoder_ABC132_A.java"
$ws.Range("D32").Value = "codeforces_373_B.java"
$ws.Range("D33").Value = "atcoder_ABC158_B.java"
$ws.Range("D34").Value = "atcoder_ABC149_B.java"
$ws.Range("D35").Value = "codeforces_669_A.java"
$ws.Range("D36").Value = "codeforces_651_A.java"
$ws.Range("D37").Value = "atcoder_ABC143_A.java"
$ws.Range("D38").Value = "codeforces_49_A.java"
$ws.Range("D39").Value = "atcoder_ABC153_A.java"
$ws.Range("D40").Value = "atcoder_ARC062_B.java"
$ws.Range("D41").Value = "atcoder_ABC170_A.java"
$ws.Range("D42").Value = "codeforces_581_A.java"
$ws.Range("D43").Value = "atcoder_ABC164_A.java"
$ws.Range("D44").Value = "atcoder_ABC174_C.java"
$ws.Range("D45").Value = "codeforces_459_A.java"
$ws.Range("D46").Value = "atcoder_AGC006_B.java"
$ws.Range("D47").Value = "atcoder_ABC070_B.java"
$ws.Range("D48").Value = "atcoder_ABC125_A.java"
$ws.Range("E48").Value = "Test Failed"
$ws.Range("D49").Value = "atcoder_AGC007_C.java"
$ws.Range("E49").Value = "Test Failed"
$ws.Range("D50").Value = "codeforces_369_B.java"
$ws.Range("E50").Value = "Test Failed"
$ws.Range("D51").Value = "codeforces_569_A.java"
$ws.Range("E51").Value = "Test Failed"
$ws.Range("D52").Value = "codeforces_30_A.java"
$ws.Range("E52").Value = "Test Failed"
$ws.Range("D53").Value = "codeforces_171_A.java"
$ws.Range("E53").Value = "Test Failed"
$ws.Range("D54").Value = "atcoder_ABC178_B.java"
$ws.Range("E54").Value = "Test Failed"
$ws.Range("D55").Value = "codeforces_190_A.java"
$ws.Range("E55").Value = "Test Failed"
$ws.Range("D56").Value = "atcoder_ABC149_C.java"
$ws.Range("E56").Value = "Test Failed"
$ws.Range("D57").Value = "codeforces_678_A.java"
$ws.Range("E57").Value = "Test Failed"
$ws.Range("D58").Value = "codeforces_276_B.java"
$ws.Range("E58").Value = "Test Failed"
$ws.Range("D59").Value = "atcoder_ABC124_A.java"
$ws.Range("E59").Value = "Test Failed"
$ws.Range("D60").Value = "codeforces_546_A.java"
$ws.Range("E60").Value = "Test Failed"
$ws.Range("D61").Value = "atcoder_AGC046_A.java"
$ws.Range("E61").Value = "Test Failed"
$ws.Range("D62").Value = "atcoder_ABC043_B.java"
$ws.Range("E62").Value = "Test Failed"
$ws.Range("D63").Value = "codeforces_622_A.java"
$ws.Range("E63").Value = "Test Failed"
$ws.Range("D64").Value = "atcoder_ABC151_A.java"
$ws.Range("E64").Value = "Test Failed"
$ws.Range("D65").Value = "atcoder_ABC158_A.java"
$ws.Range("E65").Value = "Test Failed"
$ws.Range("D66").Value = "codeforces_340_A.java"
$ws.Range("E66").Value = "Test Failed"
$ws.Range("D67").Value = "atcoder_ABC139_B.java"
$ws.Range("E67").Value = "Test Failed"
$ws.Range("D68").Value = "atcoder_ABC168_C.java"
$ws.Range("E68").Value = "Test Failed"
$ws.Range("D69").Value = "atcoder_AGC002_A.java"
$ws.Range("E69").Value = "Test Failed"
$ws.Range("D70").Value = "atcoder_AGC046_B.java"
$ws.Range("E70").Value = "Test Failed"
$ws.Range("D71").Value = "codeforces_242_A.java"
$ws.Range("E71").Value = "Test Failed"
$ws.Range("D72").Value = "codeforces_59_A.java"
$ws.Range("E72").Value = "Test Failed"
$ws.Range("D73").Value = "atcoder_ABC127_B.java"
$ws.Range("E73").Value = "Test Failed"
$ws.Range("D74").Value = "codeforces_306_A.java"
$ws.Range("E74").Value = "Test Failed"
$ws.Range("D75").Value = "atcoder_ABC142_A.java"
$ws.Range("E75").Value = "Test Failed"
$ws.Range("D76").Value = "codeforces_203_A.java"
$ws.Range("E76").Value = "Test Failed"
$ws.Range("D77").Value = "codeforces_379_A.java"
$ws.Range("E77").Value = "Test Failed"
$ws.Range("D78").Value = "atcoder_ABC051_A.java"
$ws.Range("E78").Value = "Test Failed"
$ws.Range("D79").Value = "codeforces_79_A.java"
$ws.Range("E79").Value = "Test Failed"
$ws.Range("D81").Value = "codeforces_334_A.java"
